# Fruta / hortaliza, semanal
# Insert a new weekly observation row at row 8 (Vega Modelo de Temuco - Níspero),
# pushing all the following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 8; existing rows 8..31 shift to 9..32.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with this week's observation.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "Vega Modelo de Temuco"
$ws.Range("C8").Value = "La Araucanía"
$ws.Range("D8").Value = 45243
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100104
$ws.Range("H8").Value = "Frutos de pepita"
$ws.Range("I8").Value = 100104004
$ws.Range("J8").Value = "Níspero"
$ws.Range("K8").Value = "Californiana(o)"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 22000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 22000
$ws.Range("Q8").Value = "$/bandeja 5 kilos"
$ws.Range("R8").Value = "Provincia de Quillota"
$ws.Range("S8").Value = 4400
$ws.Range("T8").Value = 5
